$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log a new volunteer-hours entry (4:32PM 1-20-2018 -> 7:23PM 1-20-2018, 171 minutes)
# in row 24, directly above the existing "Total Project Hours:" summary row (row 29).
$ws.Range("A24").Value = "4:32PM 1-20-2018"
$ws.Range("B24").Value = "7:23PM 1-20-2018"
$ws.Range("C24").Value = 171

# Move the active selection to A14
$ws.Range("A14").Select()
